# Apply the edits described by the diff to the active worksheet (Feuil1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# A3 / A11: new cells containing "coarse".
$ws.Range("A3").Value = "coarse"
$ws.Range("A11").Value = "coarse"

# A2: "T15_8" -> "T15_8 " (trailing space added) and give it the
# wrap-text style used by the header row (style index 1 == wrapText).
$ws.Range("A2").Value = "T15_8 "
$ws.Range("A2").WrapText = $true

# Move the active selection to A14 (was F16).
$ws.Range("A14").Select()
